$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.037.31"
$ws.Range("E2").Value = "  +1.49%  "
$ws.Range("D3").Value = "3.207.51"
$ws.Range("E3").Value = "  +1.05%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.81"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.51%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "3.206.00"
$ws.Range("E8").Value = "  +0.96%  "
$ws.Range("E9").Value = "  +0.16%  "
$ws.Range("E10").Value = "  -1.46%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.11"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.54%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.509"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.02%  "
$ws.Range("E13").Value = "  +0.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "39.45"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.84%  "
$ws.Range("D15").Value = "3.734.55"
$ws.Range("E15").Value = "  +1.10%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.49"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.21%  "
$ws.Range("D17").Value = "66.075.87"
$ws.Range("E17").Value = "  +1.45%  "
$ws.Range("D18").Value = "3.212.31"
$ws.Range("E18").Value = "  +1.29%  "
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "511.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.45"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.66%  "
$ws.Range("E22").Value = "  +1.76%  "
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.13"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.68%  "
$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.36"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.64%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.98"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.23%  "
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.31"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.81%  "
$ws.Range("E28").Value = "  +2.74%  "
$ws.Range("E29").Value = "  +3.36%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.89"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.89%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.85"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +8.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.12"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.82%  "
$ws.Range("E33").Value = "  +2.17%  "
$ws.Range("E34").Value = "  +0.16%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.58"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.45%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.01"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.63%  "
$ws.Range("E37").Value = "  +0.36%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "487.04"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.17%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0420"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.96"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.86%  "
$ws.Range("E41").Value = "  +2.65%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.301"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.62%  "
$ws.Range("E43").Value = "  +2.19%  "
$ws.Range("D44").Value = "2.956.62"
$ws.Range("E44").Value = "  -4.02%  "
$ws.Range("B45").Value = "PEPE"
$ws.Range("C45").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D45").Value = "0.0₃0646"
$ws.Range("E45").Value = "  +6.29%  "
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.45"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.61"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.60%  "
$ws.Range("E48").Value = "  +0.05%  "
$ws.Range("E49").Value = "  +0.96%  "
$ws.Range("E50").Value = "  +2.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "120.18"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.18%  "
